$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.267.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.945.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.08%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0872"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.406.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.934.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.983"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.306.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.04%  "

$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0984"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.181"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +20.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.108"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0445"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.84%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.63%  "

$ws.Range("E40").Value = "  -3.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("E42").Value = "  +2.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("E47").Value = "  -5.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.134.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.249"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0352"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
